# Commit: "reduce margin values to minimum for printing"
#
# The only content-bearing change in the target revision is the page
# margins on the document's (single) section -- shrunk down to the
# smallest values the printer/page allows:
#
#   top="720" right="720" bottom="720" left="720" header="708" footer="708"
#     -> top="238" right="340" bottom="249" left="340" header="709" footer="709"
#
# (all values are twentieths of a point / twips; PageSetup.* properties
# on the Word object model are expressed in points, so divide by 20).

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    $ps = $sec.PageSetup
    $ps.TopMargin = 238 / 20.0
    $ps.RightMargin = 340 / 20.0
    $ps.BottomMargin = 249 / 20.0
    $ps.LeftMargin = 340 / 20.0
    $ps.HeaderDistance = 709 / 20.0
    $ps.FooterDistance = 709 / 20.0
}
